# Update calibration estimates (std moving average) across the
# riskDriverType PnL factor calibration sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------
$wsAR = $wb.Worksheets.Item("AR")
$wsAR.Range("B2").Value = -0.0009957742384762028
$wsAR.Range("B3").Value = 0.7808667006718862
$wsAR.Range("B4").Value = 0.1062094442685016
$wsAR.Range("B5").Value = "[0.9999999999999999, 0.02562816831403599, -0.06473359358270972, -0.08551161735630672, -0.08531476330581801, 0.0863416177961257, 0.022616504208740573, 0.03853118877700369, 0.04587237568491471, -0.002139168877255181, -0.02079972463942587, -0.05021462530082704, -0.05733000882239452, 0.013532331744224281, 0.08611254494226606, 0.0028074846081215885, -0.011158655111631069, -0.047232466496148366, -0.02468375501175664, -0.01739956055248122]"

# ---------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------
$wsSETAR = $wb.Worksheets.Item("SETAR")
$wsSETAR.Range("B4").Value = -0.1901017084119081
$wsSETAR.Range("B5").Value = 0.5620286220926134
$wsSETAR.Range("B6").Value = 0.07157351391273913
$wsSETAR.Range("B7").Value = 0.1606608288492219
$wsSETAR.Range("B8").Value = 0.5982193211705331
$wsSETAR.Range("B9").Value = 0.05976853915506266
$wsSETAR.Range("B10").Value = "[1.0, 0.07056668287085331, -0.020067523868945572, -0.023510646431827376, 0.009313353938934094, 0.04011340161391695, -0.006836439378172958, 0.00935751665841219, -0.001911791425266478, -0.024539884366682863, -0.004839997310232216, -0.009481277301696569, -0.02229988103688684, -0.02292065599292571, 0.019110836761371052, 0.004348848697662681, -0.005205494702985697, -0.04836938060319221, -0.004511300545584137, 0.01898771571660957]"

# ---------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------
$wsGARCH = $wb.Worksheets.Item("GARCH")
$wsGARCH.Range("B2").Value = -0.0004018072948923263
$wsGARCH.Range("B3").Value = 0.1153235563841152
$wsGARCH.Range("B4").Value = 0.03270380227550945
$wsGARCH.Range("B5").Value = 0.0000003427989374063351
$wsGARCH.Range("B6").Value = "[1.0, 0.011247879761085263, -0.06299864758850016, -0.07190144356188069, -0.10338243902158446, 0.1004769435354579, 0.013697820983429658, 0.0287948861368788, 0.06088941770048124, 0.015836162447387184, -0.055893381317907685, -0.05125953952864348, -0.05724729631573349, 0.022141167853776227, 0.08781888311757843, 0.011550131154847718, -0.016252455264519038, -0.03919928374510832, -0.05056645656772844, -0.03775074619122891]"

# ---------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------
$wsTARCH = $wb.Worksheets.Item("TARCH")
$wsTARCH.Range("B2").Value = -0.001694176926069925
$wsTARCH.Range("B3").Value = 0.1156572388749257
$wsTARCH.Range("B4").Value = 0.00000000003726031878575412
$wsTARCH.Range("B5").Value = 0
$wsTARCH.Range("B6").Value = "[1.0, 0.011679048482686551, -0.06036469505124449, -0.0717250945731754, -0.10340309290650686, 0.1007683208413346, 0.01215136090052689, 0.027800050377162162, 0.06129954739383657, 0.014602282831571134, -0.05547950722722275, -0.05143517987312567, -0.0556101731684857, 0.021731271506579483, 0.08679080188171523, 0.012559447323146723, -0.015969915741991924, -0.0404104804423034, -0.05232143061860133, -0.036894588088003825]"
$wsTARCH.Range("B7").Value = 0.05810480409792006

# ---------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------
$wsARTARCH = $wb.Worksheets.Item("AR_TARCH")
$wsARTARCH.Range("B2").Value = -0.001007300465803767
$wsARTARCH.Range("B3").Value = 0.1036995315602178
$wsARTARCH.Range("B4").Value = 0.00000000002054371562639695
$wsARTARCH.Range("B5").Value = 0
$wsARTARCH.Range("B6").Value = "[1.0, 0.0011163013921449026, -0.06361130523716504, -0.08430205860625888, -0.08794262397825933, 0.0895346808390797, 0.018577849378341994, 0.036797548140364726, 0.04558582266024879, -0.0038070266527379594, -0.02044937331667911, -0.04941918262412142, -0.056335536675284625, 0.012904348628622735, 0.08605150220830564, 0.0025501833924675096, -0.009760708297464293, -0.047258470573420426, -0.0261954886188416, -0.01653902268542312]"
$wsARTARCH.Range("B7").Value = 0.04483924486822021
$wsARTARCH.Range("B9").Value = 0.7813370938224929
